$d = $word.ActiveDocument

# 1. Update activation date
$d.Content.Find.Execute("Ativação: 01/01/2016", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2022", 2) | Out-Null

# 2. Replace Objetivos paragraph text
$d.Content.Find.Execute("Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar. Apresentar como se classificam os materiais em função de suas propriedades. Aspectos legais e o mercado do engenheiro de materiais.", $true, $false, $false, $false, $false, $true, 1, $false, "A disciplina busca introduzir o aluno ao ambiente de engenharia, propondo problemas desafiadores gerando aptidão para solução de problemas. Apresentar a Engenharia de Materiais e seus campos de atuação, aspectos legais e éticos, bem como o mercado de trabalho para o engenheiro de materiais no Século XXI. Propiciar aos alunos uma visão geral do curso, com apresentação do currículo do curso de Engenharia de Materiais da EEL. Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar.", 2) | Out-Null

# 3. Add second docente responsavel with line break (split into two runs,
#    matching the source structure) by rewriting the whole paragraph's XML.
$found = $d.Content
$found.Find.Execute("984972 - Hugo Ricardo Zschommler Sandim") | Out-Null
$p = $found.Paragraphs.First
$paraRng = $d.Range($p.Range.Start, $p.Range.End + 1)
$docenteXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>984972 - Hugo Ricardo Zschommler Sandim</w:t><w:br/></w:r><w:r><w:t>7459752 - Maria Ismenia Sodero Toledo Faria</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraRng.InsertXML($docenteXml)

# 4. Replace Programa resumido text
$d.Content.Find.Execute("A importância dos materiais na evolução do homem na pré-história. Alquimia e Revolução Científica. Materiais no século XIX e a Revolução Industrial. Revolução do século XX e os materiais nas guerras mundiais.A Ciência e Engenharia de Materiais como Interdisciplina. Ciclo dos materiais. O Tetraedro da Ciência e Engenharia de Materiais. Classificação dos materiais em função de suas propriedades. Noções de seleção de materiais. Diagramas de Ashby. Estudos de casos. Legislação. O mercado de trabalho para o engenheiro de materiais. Perspectivas para a Ciência e Engenharia de Materiais. Visita técnica a uma grande empresa do setor.", $true, $false, $false, $false, $false, $true, 1, $false, "1- A importância dos materiais na evolução do homem na pré-história. Alquimia, Revolução Científica e a Revolução Industrial. 2-O Engenheiro como um profissional, funções da engenharia, a ética e comunicação na engenharia 3-A grandes áreas da Engenharia de Materiais. A interdisciplinaridade da Ciência e Engenharia de Materiais. 4- Perspectivas para a Engenharia de Materiais no século XXI. 5- O currículo do curso de engenharia de materiais da EEL-USP. 6- Noções básicas de Projetos em Engenharia.Em todos o conteúdo do curso serão abordados aspectos sociais, ambientais, éticos, legais e econômicos para ampliar as competências dos alunos", 2) | Out-Null

# 5. Replace Programa text
$d.Content.Find.Execute("Parte 1 - Materiais na história e na sociedade.1) A importância dos materiais na evolução do homem na pré-história.2) As idades do cobre, do bronze e do ferro. Materiais empregados nas eras clássica e medieval.3) Alquimia e Revolução Científica.4) Materiais no século XIX e a Revolução Industrial.5) Revolução do século XX e os materiais nas guerras mundiais.Parte 2 - A Ciência e Engenharia de Materiais como Interdisciplina.6) Ciclo dos materiais. O Tetraedro da Ciência e Engenharia de Materiais. 7) Classificação dos materiais em função de suas propriedades. Noções de seleção de materiais. Diagramas de Ashby.8) Estudos de casos.9) Legislação. O mercado de trabalho para o engenheiro de materiais.10) Perspectivas para a Ciência e Engenharia de Materiais.Conteúdo prático: 1. Visita ao Departamento de Engenharia de Materiais. Visita externa para integralização dos conhecimentos.", $true, $false, $false, $false, $false, $true, 1, $false, "1- As características importantes de um engenheiro: aptidões interpessoais, aptidões de comunicação, liderança e competência. O engenheiro, profissional que busca solucionar problemas. 2-A Engenharia de Materiais: áreas de atuação e mercado de trabalho. Aplicação. A importância dos materiais na evolução do homem, as grandes áreas e interdisciplinaridade da Ciência e Engenharia de Materiais. Visita ao Departamento de Engenharia de Materiais. Conhecimento dos Grupos de Pesquisa do Departamento. Perspectivas para a Engenharia de Materiais no século XXI. 3- O campo de trabalho do engenheiro de materiais e suas áreas de atuação. Visita externa para integralização dos conhecimentos. 4- O currículo do curso de engenharia de materiais na EEL/USP. 5- Apresentação do método de trabalho com projetos, definindo os atributos de um projeto de engenharia, mapas conceituais e ferramentas que ilustram ideias e relações entre elas. Formular estratégias para resolução de problemas de engenharia. Estudo de casos", 2) | Out-Null

# 6. Replace Metodo text
$d.Content.Find.Execute("O aluno será avaliado ao longo do semestre por duas avaliações escritas (P1 e P2) e com pesos iguais.", $true, $false, $false, $false, $false, $true, 1, $false, "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras", 2) | Out-Null

# 7. Replace Criterio text
$d.Content.Find.Execute("Nota Final NF = [P1 + P2]/2", $true, $false, $false, $false, $false, $true, 1, $false, "Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas.", 2) | Out-Null

# 8. Replace Norma de recuperacao text
$d.Content.Find.Execute("Para a recuperação será realizada uma prova escrita (PR) abrangendo toda a matéria lecionada no semestre, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2) | Out-Null

# 9. Replace Bibliografia text
$d.Content.Find.Execute("1) Cohem, M. Ciência e Engenharia de Materias: Sua Evolução, Prática e Perspectivas.Parte I - Materiais na História e na Sociedade.Parte II - A Ciência e Engenharia de Materiais como uma Multidisciplina.Tradução José Roberto da Silva. DEMa/UFSCar - São Carlos - 1981. 2) Callister Jr., W.D.C. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 5a.ed., 2002.3) Understanding Materials Science: History, Properties, Applications.  ROLF E. HUMMEL.  Springer, 1997.4) Sustainable Development and the Advanced Materials: The Brazilian Case. Ed. Roberto C. Vilas Boas . CETEM -  Centro de Tecnologia Mineral MCT/CNPq, IDRC/Canadá - International Development Research Center, 1995.", $true, $false, $false, $false, $false, $true, 1, $false, "1) BROCKMAN, J.B. Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2010.2) M.T. HOLTZAPPLE, W.D. REECE, Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2006.2) CALLISTER Jr., W.D. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 7a.ed., 2008. 4) - COHEN, M. (Ed.). Ciência e Engenharia de Materiais: sua Evolução, Prática e Perspectivas. Parte I: Materiais na história e na sociedade, 98p. Parte II: A Ciência e Engenharia de Materiais como uma multidisciplina, Tradução: José Roberto Gonçalves da Silva, São Carlos, UFSCar, 1985.5) Artigos científicos", 2) | Out-Null

